# Generate Report for Handoff
# File "b.md" moves from "Handed back: in sync with en-US" to "Ready for handoff"
# and records new handoff file / datetime for both zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: update status for b.md (row 3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: update status / latest handoff file / latest handoff datetime
# for b.md (row 3)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-10 22:38:10"

# The hyperlink display text for C3 needs to match the new file name while
# keeping the same target URL. Rebuild all hyperlinks on the sheet (in order)
# so relationship ids / targets stay identical except for the updated display
# text on C3.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a4361e4a0be6c35e974f2ac828933fd689315d5a/e2e/a.md", "", "", "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a1df321053f6ceb6c1a94fd3a13c7711c2a7ad0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b437363316e791b63cfeb58a7aecc69c24f4a1f7/e2e/a.md", "", "", "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/99a1de1804f3133d4d0cfece7aa7813d6d22abb9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a4361e4a0be6c35e974f2ac828933fd689315d5a/e2e/b.md", "", "", "b.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a1df321053f6ceb6c1a94fd3a13c7711c2a7ad0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b437363316e791b63cfeb58a7aecc69c24f4a1f7/e2e/a.md", "", "", "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/99a1de1804f3133d4d0cfece7aa7813d6d22abb9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a4361e4a0be6c35e974f2ac828933fd689315d5a/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: update status / latest handoff file / latest handoff datetime
# for b.md (row 3)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-10 22:38:19"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a4361e4a0be6c35e974f2ac828933fd689315d5a/e2e/a.md", "", "", "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/975568a5d3dda9d876bc6e6292b78ee996dcd08b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/782ec2060511afc26327dba4d745221db5a38a17/e2e/a.md", "", "", "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/717aefa928898ec64a4484c5ac0a51d9fec1e955/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a4361e4a0be6c35e974f2ac828933fd689315d5a/e2e/b.md", "", "", "b.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/975568a5d3dda9d876bc6e6292b78ee996dcd08b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/782ec2060511afc26327dba4d745221db5a38a17/e2e/a.md", "", "", "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/717aefa928898ec64a4484c5ac0a51d9fec1e955/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a4361e4a0be6c35e974f2ac828933fd689315d5a/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Done applying handoff report update."
